$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Value = "[{`"timestamp`":`"30 September 2020, 8.52`",`"keyCount`":318,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"L\/VJRgfcQoXAEqA7YqpW4JjXXUm4qTk2ETXOxXI\/C0A=`"},{`"timestamp`":`"29 September 2020, 8.09`",`"keyCount`":238,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"buIXTow1ZHDIVMQhEmgOqDBqbydn03Gm8NGhy6i41Tc=`"},{`"timestamp`":`"28 September 2020, 8.07`",`"keyCount`":199,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"goyTayM0mYovHuxJ0OSDqmm+oB3bTRmYvy9uWbmVrjQ=`"},{`"timestamp`":`"27 September 2020, 8.05`",`"keyCount`":141,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"bfObmmKf4yYCy3kecvhRL4d57FZ8UxUfFlrfWq7YzLg=`"},{`"timestamp`":`"26 September 2020, 8.01`",`"keyCount`":242,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"AP5aEpXLoQlr0QmfeoIh4FHy01c+H5xDR4JCbPprdPc=`"},{`"timestamp`":`"25 September 2020, 8.00`",`"keyCount`":217,`"matchesCount`":1,`"appName`":`"Koronavilkku`",`"hash`":`"nLBOmHqB7xON3sOJLjm86EVtaADO9z8MqVcA8r3v4t8=`"},{`"timestamp`":`"24 September 2020, 7.19`",`"keyCount`":211,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"XlRtacwf5zNbpVVHoSKu4IDF9qs1Wh6PipZr1sPMn+E=`"},{`"timestamp`":`"23 September 2020, 7.15`",`"keyCount`":189,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"XraKuKwjDha2BdTzS3+Xm08W+39n7AwLVWCl6VKlkHM=`"},{`"timestamp`":`"22 September 2020, 7.13`",`"keyCount`":311,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"cg00ZNzVnYnP9YAuDqm5ML6GrpkFfkESOfEZPrL6R1c=`"},{`"timestamp`":`"21 September 2020, 7.11`",`"keyCount`":157,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"eiScOGMBbruz9+mpIa4JU4R5C1NlJKl+7o3xWmhO7jw=`"},{`"timestamp`":`"20 September 2020, 11.09`",`"keyCount`":202,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"AT2kruWiYpVjiscFZlBLAG6pFxBrAq3cpo0ZIHpx2uM=`"},{`"timestamp`":`"19 September 2020, 10.50`",`"keyCount`":190,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"wQ25zYYr\/T3LDtbukS3OYyUSNBYzmsEMqI3tDEx+GdY=`"},{`"timestamp`":`"18 September 2020, 10.37`",`"keyCount`":82,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"LxnzOI42gpmBHhNssMRQgJbJyyVz6j95FgbhcEhslFE=`"},{`"timestamp`":`"17 September 2020, 10.36`",`"keyCount`":137,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"nPOu6VJOVmdCS0JgyedgK1KtW4rxnqzJUGhlNwGuOV4=`"},{`"timestamp`":`"16 September 2020, 10.32`",`"keyCount`":125,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"L5x\/yuwZx1S3i2n0\/9FRe4HclDlgucBDfy5Or56ZchQ=`"},{`"timestamp`":`"15 September 2020, 10.28`",`"keyCount`":136,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"op7f5yD\/eSrwkpUHdNNBhVmT48knZIc26Jzt8vGL3JA=`"},{`"timestamp`":`"14 September 2020, 10.26`",`"keyCount`":67,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"vIEKsogtvkNLwDVfjkgES\/2vVVF0c16tzUm34F9Dse4=`"},{`"timestamp`":`"13 September 2020, 10.22`",`"keyCount`":87,`"matchesCount`":0,`"appName`":`"Koronavilkku`",`"hash`":`"361cWfGKKo4KQRNrZy3CtjvAmySVn72Z3X1Sx903ClE=`"}]"
